# Update "doctor_MA average" (column AF) figures in rows 4-13 with the
# newly computed results, per commit "updated results and code".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AF4").Value  = 0.583
$ws.Range("AF5").Value  = 0.857
$ws.Range("AF6").Value  = 0.694
$ws.Range("AF7").Value  = 0.783
$ws.Range("AF8").Value  = 0.804
$ws.Range("AF9").Value  = 0.714
$ws.Range("AF10").Value = 0.857
$ws.Range("AF11").Value = 0.857
$ws.Range("AF12").Value = 1.167
$ws.Range("AF13").Value = 2
